# labs: finalise lab 02
# Duplicate the "SO4" worksheet, name the copy "SO4_excel", place it as the
# last (3rd) sheet in the workbook, and make it the active/selected tab.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("SO4")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy "SO4" to a new sheet positioned immediately after the current last sheet.
$source.Copy($null, $lastSheet)

# The freshly-copied sheet is now the last sheet in the workbook.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "SO4_excel"

# Make the new sheet the active tab, matching the saved workbook view.
$newSheet.Activate()
